$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 161, shifting existing rows 161-200 down to 162-201.
$ws.Rows.Item(161).Insert()

# Populate the new row 161 with the new weekly price record.
$ws.Cells.Item(161, 1).Value = 5
$ws.Cells.Item(161, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(161, 3).Value = "Maule"
$ws.Cells.Item(161, 4).Value = 44508
$ws.Cells.Item(161, 4).NumberFormat = $ws.Cells.Item(162, 4).NumberFormat
$ws.Cells.Item(161, 5).Value = 7
$ws.Cells.Item(161, 6).Value = 100114014
$ws.Cells.Item(161, 7).Value = "Betarraga"
$ws.Cells.Item(161, 8).Value = "Sin especificar"
$ws.Cells.Item(161, 9).Value = "Primera"
$ws.Cells.Item(161, 10).Value = 5000
$ws.Cells.Item(161, 11).Value = 650
$ws.Cells.Item(161, 12).Value = 650
$ws.Cells.Item(161, 13).Value = 650
$ws.Cells.Item(161, 14).Value = "$/paquete 5 unidades"
$ws.Cells.Item(161, 15).Value = "Región del Maule"
$ws.Cells.Item(161, 16).Value = 130
$ws.Cells.Item(161, 17).Value = 5
$ws.Cells.Item(161, 18).Value = "Hortaliza"
